$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:E").Insert()
$addr = $ws.UsedRange.Address()
Write-Host "Dimension after insert: $addr"
$d = $ws.Cells.Item(7,4).Value2
Write-Host "Cell(7,4) after insert: $d"
$f = $ws.Cells.Item(7,6).Value2
Write-Host "Cell(7,6) after insert: $f"
